# Unified file classes in the class page
#
# The "class" page lists the number of search results found. Previously
# only the "7 results have been found." message existed on the sheet;
# this adds the other result-count messages (singular "1 result" and
# plural "4 results") next to the existing rows, as cell B2 and B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "1 result has been found."
$ws.Range("B3").Value = "4 results have been found."
